$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# 1. Change the weekly eating-out budget formula from 25 to 20 (24.5*20*4)
$ws1.Range("B5").Formula = "=24.5*20*4"

# 2. Update rent: add a $50 cell to the base 1735, making rent 1785
$ws1.Range("F9").Value = 50
$ws1.Range("E9").Formula = "=1735+F9"

# 3. Update insurance/owed amount in H10 from 120 to 130
$ws1.Range("H10").Value = 130

# 4. Update the annual savings label to reflect new total (9600 saved per year)
$ws1.Range("A16").Value = "9600 saved per year"

# 5. Rename "chipolte" expense label to "Eating out"
$ws1.Range("D4").Value = "Eating out"

# Update selections / active sheet/tab to match final state
[void]$ws1.Range("E9").Select()
[void]$ws2.Range("B9").Select()
[void]$ws2.Activate()
